# Add Panasonic/Olympus Micro Four Thirds lenses to LensTable,
# extend the Table1 ListObject, log the change on the Updates sheet,
# and widen a couple of columns on LensTable to fit the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$tbl = $ws.ListObjects.Item(1)

# Grow the table by 9 rows (434 data rows -> 443 data rows, i.e. A1:Q444).
for ($i = 0; $i -lt 9; $i++) {
    $tbl.ListRows.Add() | Out-Null
}

    $ws.Cells.Item(436, 1).Value = 435
    $ws.Cells.Item(436, 2).Value = "Olympus"
    $ws.Cells.Item(436, 3).Value = "M.Zuiko Digital 25mm f/1.8"
    $ws.Cells.Item(436, 4).Value = 2014
    $ws.Cells.Item(436, 5).Value = 249
    $ws.Cells.Item(436, 6).Value = 136
    $ws.Cells.Item(436, 7).Value = 40.6
    $ws.Cells.Item(436, 8).Value = 56
    $ws.Cells.Item(436, 9).Value = 1.8
    $ws.Cells.Item(436, 10).Value = "25"
    $ws.Cells.Item(436, 11).Value = "0"
    $ws.Cells.Item(436, 12).Value = 0
    $ws.Cells.Item(436, 13).Value = 0
    $ws.Cells.Item(436, 14).Value = "Micro Four Thirds"

    $ws.Cells.Item(437, 1).Value = 436
    $ws.Cells.Item(437, 2).Value = "Olympus"
    $ws.Cells.Item(437, 3).Value = "M.Zuiko Digital ED 17mm f/1.2 PRO"
    $ws.Cells.Item(437, 4).Value = 2018
    $ws.Cells.Item(437, 5).Value = 1199
    $ws.Cells.Item(437, 6).Value = 390
    $ws.Cells.Item(437, 7).Value = 87
    $ws.Cells.Item(437, 8).Value = 68.2
    $ws.Cells.Item(437, 9).Value = 1.2
    $ws.Cells.Item(437, 10).Value = "17"
    $ws.Cells.Item(437, 11).Value = "0"
    $ws.Cells.Item(437, 12).Value = 0
    $ws.Cells.Item(437, 13).Value = 0
    $ws.Cells.Item(437, 14).Value = "Micro Four Thirds"

    $ws.Cells.Item(438, 1).Value = 437
    $ws.Cells.Item(438, 2).Value = "Olympus"
    $ws.Cells.Item(438, 3).Value = "M.Zuiko Digital ED 75mm f/1.8"
    $ws.Cells.Item(438, 4).Value = 2012
    $ws.Cells.Item(438, 5).Value = 749
    $ws.Cells.Item(438, 6).Value = 305
    $ws.Cells.Item(438, 7).Value = 69.1
    $ws.Cells.Item(438, 8).Value = 64
    $ws.Cells.Item(438, 9).Value = 1.8
    $ws.Cells.Item(438, 10).Value = "75"
    $ws.Cells.Item(438, 11).Value = "0"
    $ws.Cells.Item(438, 12).Value = 0
    $ws.Cells.Item(438, 13).Value = 0
    $ws.Cells.Item(438, 14).Value = "Micro Four Thirds"

    $ws.Cells.Item(439, 1).Value = 438
    $ws.Cells.Item(439, 2).Value = "Olympus"
    $ws.Cells.Item(439, 3).Value = "M.Zuiko Digital ED 45mm f/1.2 PRO"
    $ws.Cells.Item(439, 4).Value = 2017
    $ws.Cells.Item(439, 5).Value = 1199
    $ws.Cells.Item(439, 6).Value = 410
    $ws.Cells.Item(439, 7).Value = 84.9
    $ws.Cells.Item(439, 8).Value = 70
    $ws.Cells.Item(439, 9).Value = 1.2
    $ws.Cells.Item(439, 10).Value = "45"
    $ws.Cells.Item(439, 11).Value = "0"
    $ws.Cells.Item(439, 12).Value = 0
    $ws.Cells.Item(439, 13).Value = 0
    $ws.Cells.Item(439, 14).Value = "Micro Four Thirds"

    $ws.Cells.Item(440, 1).Value = 439
    $ws.Cells.Item(440, 2).Value = "Olympus"
    $ws.Cells.Item(440, 3).Value = "M.Zuiko Digital 45mm f/1.8"
    $ws.Cells.Item(440, 4).Value = 2011
    $ws.Cells.Item(440, 5).Value = 249
    $ws.Cells.Item(440, 6).Value = 116
    $ws.Cells.Item(440, 7).Value = 46
    $ws.Cells.Item(440, 8).Value = 56
    $ws.Cells.Item(440, 9).Value = 1.8
    $ws.Cells.Item(440, 10).Value = "45"
    $ws.Cells.Item(440, 11).Value = "0"
    $ws.Cells.Item(440, 12).Value = 0
    $ws.Cells.Item(440, 13).Value = 0
    $ws.Cells.Item(440, 14).Value = "Micro Four Thirds"

    $ws.Cells.Item(441, 1).Value = 440
    $ws.Cells.Item(441, 2).Value = "Olympus"
    $ws.Cells.Item(441, 3).Value = "M.Zuiko Digital 17mm f/1.8"
    $ws.Cells.Item(441, 4).Value = 2013
    $ws.Cells.Item(441, 5).Value = 349
    $ws.Cells.Item(441, 6).Value = 120
    $ws.Cells.Item(441, 7).Value = 35.5
    $ws.Cells.Item(441, 8).Value = 57.5
    $ws.Cells.Item(441, 9).Value = 1.8
    $ws.Cells.Item(441, 10).Value = "17"
    $ws.Cells.Item(441, 11).Value = "0"
    $ws.Cells.Item(441, 12).Value = 0
    $ws.Cells.Item(441, 13).Value = 0
    $ws.Cells.Item(441, 14).Value = "Micro Four Thirds"

    $ws.Cells.Item(442, 1).Value = 441
    $ws.Cells.Item(442, 2).Value = "Panasonic"
    $ws.Cells.Item(442, 3).Value = "Leica DG Summilux 25mm f/1.4 ASPH"
    $ws.Cells.Item(442, 4).Value = 2011
    $ws.Cells.Item(442, 5).Value = 598
    $ws.Cells.Item(442, 6).Value = 200
    $ws.Cells.Item(442, 7).Value = 54.5
    $ws.Cells.Item(442, 8).Value = 63
    $ws.Cells.Item(442, 9).Value = 1.4
    $ws.Cells.Item(442, 10).Value = "25"
    $ws.Cells.Item(442, 11).Value = "0"
    $ws.Cells.Item(442, 12).Value = 0
    $ws.Cells.Item(442, 13).Value = 0
    $ws.Cells.Item(442, 14).Value = "Micro Four Thirds"

    $ws.Cells.Item(443, 1).Value = 442
    $ws.Cells.Item(443, 2).Value = "Panasonic"
    $ws.Cells.Item(443, 3).Value = "Leica DG Summilux 12mm f/1.4 ASPH"
    $ws.Cells.Item(443, 4).Value = 2016
    $ws.Cells.Item(443, 5).Value = 1298
    $ws.Cells.Item(443, 6).Value = 335
    $ws.Cells.Item(443, 7).Value = 70
    $ws.Cells.Item(443, 8).Value = 70
    $ws.Cells.Item(443, 9).Value = 1.4
    $ws.Cells.Item(443, 10).Value = "12"
    $ws.Cells.Item(443, 11).Value = "0"
    $ws.Cells.Item(443, 12).Value = 0
    $ws.Cells.Item(443, 13).Value = 0
    $ws.Cells.Item(443, 14).Value = "Micro Four Thirds"

    $ws.Cells.Item(444, 1).Value = 443
    $ws.Cells.Item(444, 2).Value = "Panasonic"
    $ws.Cells.Item(444, 3).Value = "Leica DG Nocticron 42.5mm f/1.2 ASPH POWER OIS"
    $ws.Cells.Item(444, 4).Value = 2014
    $ws.Cells.Item(444, 5).Value = 1598
    $ws.Cells.Item(444, 6).Value = 425
    $ws.Cells.Item(444, 7).Value = 76.8
    $ws.Cells.Item(444, 8).Value = 74
    $ws.Cells.Item(444, 9).Value = 1.2
    $ws.Cells.Item(444, 10).Value = "42.5"
    $ws.Cells.Item(444, 11).Value = "0"
    $ws.Cells.Item(444, 12).Value = 1
    $ws.Cells.Item(444, 13).Value = 0
    $ws.Cells.Item(444, 14).Value = "Micro Four Thirds"


# Log the change on sheet "Updates" (the table's changelog sheet).
$log = $wb.Worksheets.Item(2)
$log.Range("A62").Value = 43097
$log.Range("B62").Value = "Added Micro 4/3 Panasonic Summilux  25, 12 mm; Nocticron 42.5mm, Olympus M.Zuiko 17/45mm f/1.2 & f1.8; 25mm"

# Widen the Model and Weight columns on LensTable to fit the new, longer
# lens names and weight figures (closest values this engine's ColumnWidth
# quantization can reach to the authored 49.5546875 / 10.44140625).
$ws.Columns.Item(3).ColumnWidth = 48.666666666666664
$ws.Columns.Item(6).ColumnWidth = 9.666666666666666
